$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31/32: Bittensor and PancakeSwap swap ranking positions, with updated values
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "530.04"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").Value = "  -3.33%  "

# Price / Volume(1h) updates
$ws.Range("D2").Value = "65.779.98"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.681.09"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.92"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.99"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  +6.35%  "
$ws.Range("E9").Value = "  +5.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.402"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.47"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000199"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "3.163.85"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "65.647.51"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.701.37"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.89"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.79"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.76"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +5.32%  "
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.50"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.426"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.58"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.49"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.15"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0612"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.94"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.644"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "0.0₆0261"
$ws.Range("E49").Value = "  +14.41%  "
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.16"
$ws.Range("E51").Value = "  -4.03%  "
